$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($i = 1; $i -le 13; $i++) {
    $row = $i + 1
    $ws.Range("A$row").Value = "s$i"
}

$ws.Range("A15").Select()
